$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E12" = 18.04990000000002
    "E32" = 16.79999999999998
    "E36" = 16.98280000000001
    "E38" = 16.46179999999999
    "E46" = 17.05
    "E54" = 16.6471
    "E55" = 16.60830000000001
    "E67" = 17.23020000000002
    "E69" = 17.30680000000002
    "E72" = 16.47440000000001
    "E91" = 18.61510000000001
    "E99" = 16.4637
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
